$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rsnn materials ("thesis-body" and "thesis-slides" rows), shifting
# the remaining rows ("cv", "sop") up.
$ws.Rows("2:3").Delete()

# Keep the _FilterDatabase defined name in sync with the now-smaller data range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase") {
        $n.RefersTo = "=docs!`$D`$1:`$D`$4"
    }
}

$ws.Range("A2").Select()
